$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill previously empty cells (leading apostrophe forces text so leading zeros survive)
$ws.Range("C2").Value = "'016015"
$ws.Range("F2").Value = "通讯行业"

# Swap DATE_TYPE_CODE value (leading apostrophe forces text so leading zeros survive)
$ws.Range("J2").Value = "'002"

# Update report date
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Update numeric figures
$ws.Range("O2").Value = 10470879.33
$ws.Range("P2").Value = 42.5208052188
$ws.Range("Q2").Value = 68325405.84999999
$ws.Range("R2").Value = 277.4601045511
$ws.Range("S2").Value = 27483975.38
$ws.Range("T2").Value = 111.6086554854
$ws.Range("U2").Value = 356396.94
$ws.Range("V2").Value = 1.4472791051
$ws.Range("Y2").Value = 1209109.91
$ws.Range("Z2").Value = 4.9100295544
$ws.Range("AA2").Value = 13695843.72
$ws.Range("AB2").Value = 55.6169434077
$ws.Range("AC2").Value = 24625308.19

# AD2 becomes blank (was a number before)
$ws.Range("AD2").ClearContents()
